$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values per row, preserved as text
# via a leading apostrophe so Excel does not reinterpret them as numbers/percentages.
$updates = @(
    @{ Row = 2; D = "326.14"; E = "-1.23%" }
    @{ Row = 3; D = "44.94"; E = "2.80%" }
    @{ Row = 4; D = "5.556"; E = "-5.74%" }
    @{ Row = 5; D = "0.08079"; E = "-2.72%" }
    @{ Row = 6; D = "8.712"; E = "-0.96%" }
    @{ Row = 7; D = "4.331"; E = "-3.85%" }
    @{ Row = 8; D = "1.902"; E = "-2.61%" }
    @{ Row = 9; D = "2.742"; E = "-5.54%" }
    @{ Row = 10; D = "0.9483"; E = "2.12%" }
    @{ Row = 11; D = "0.1176"; E = "-6.25%" }
    @{ Row = 12; D = "0.1893"; E = "-3.04%" }
    @{ Row = 13; D = "0.1013"; E = "7.03%" }
    @{ Row = 14; D = "0.04192"; E = "5.35%" }
    @{ Row = 15; D = "0.1065"; E = "0.22%" }
    @{ Row = 16; D = "0.001281"; E = "-1.67%" }
    @{ Row = 17; D = "0.005982"; E = "-0.51%" }
    @{ Row = 18; D = "3.605"; E = "2.34%" }
    @{ Row = 19; D = "0.3486"; E = "-0.67%" }
    @{ Row = 20; D = "8.690"; E = "-4.53%" }
    @{ Row = 21; D = "0.1372"; E = "0.01%" }
    @{ Row = 22; D = "0.2662"; E = "3.44%" }
    @{ Row = 23; D = "0.04254"; E = "-3.17%" }
    @{ Row = 24; D = "0.001236"; E = "-1.63%" }
    @{ Row = 25; D = "0.004633"; E = "5.38%" }
    @{ Row = 26; D = $null; E = "3.49%" }
    @{ Row = 27; D = "0.0004000"; E = "0.13%" }
    @{ Row = 39; D = "0.02658"; E = "-4.85%" }
    @{ Row = 40; D = "0.05554"; E = "-0.82%" }
    @{ Row = 41; D = $null; E = "24.68%" }
    @{ Row = 42; D = "0.007685"; E = "-2.86%" }
    @{ Row = 43; D = "0.1394"; E = "-1.82%" }
    @{ Row = 44; D = "0.002060"; E = "-2.07%" }
    @{ Row = 45; D = "0.009220"; E = "-11.41%" }
    @{ Row = 46; D = "0.00007114"; E = "-1.41%" }
    @{ Row = 47; D = $null; E = "0.12%" }
    @{ Row = 48; D = "0.003444"; E = "-10.61%" }
    @{ Row = 49; D = $null; E = "-0.21%" }
    @{ Row = 50; D = "0.00002106"; E = "0.12%" }
    @{ Row = 51; D = "0.0002006"; E = "0.12%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E
    }
}
